$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New "Rekening" classes readied for review (rows 15-18) ---
# Rows 15 & 16 keep the plain style (s=2) already present on the template.
$ws.Range("A15").Value = "RekeningDAO"
$ws.Range("B15").Value = "Ju-Sen"

$ws.Range("A16").Value = "JDBCRekeningDAO"
$ws.Range("B16").Value = "Ju-Sen"

# Rows 17 & 18 get the highlighted style used by A13/A14 (style index 5) --
# copy that formatting across before writing the values.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "Model Rekening"
$ws.Range("B17").Value = "Ju-Sen"

$ws.Range("A13").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "RekeningService"
$ws.Range("B18").Value = "Ju-Sen"

$excel.CutCopyMode = 0

# --- Move the active selection from C10 to C16 ---
$ws.Activate() | Out-Null
$ws.Range("C16").Select() | Out-Null

Write-Host "done"
